# templateInforme.docx fixes:
#   - "arreglo de formato de tabla de vulnarabilidades" + "refactoring de
#     tabla de ips-puertos"
#
# Net effect inside this template document:
#   1) The "_GoBack" bookmark (Word's automatic "last edit position"
#      marker) now sits at the very start of the document -- right after
#      the first paragraph's pPr, before its run -- instead of sitting
#      alone in the last (empty) paragraph before the sectPr.
#   2) The built-in "Heading 1" ("Titulo 1") paragraph style is changed to
#      left justification, which is what fixes the heading alignment used
#      by the vulnerability table.

$d = $word.ActiveDocument

# --- 1) Relocate the "_GoBack" bookmark -------------------------------

# Remove it from wherever Word currently has it (the trailing empty
# paragraph of the document body).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create it spanning the first paragraph, so it ends up as a pair of
# empty bookmarkStart/bookmarkEnd markers right after that paragraph's
# pPr and before its run -- matching where Word leaves "_GoBack" after
# the most recent edit was made at the top of the document. Re-deriving
# the Range from plain numeric offsets (rather than reusing the
# paragraph's own Range object) keeps the bookmark pair anchored inside
# paragraph 1 instead of spilling its end marker into paragraph 2.
$firstPara = $d.Paragraphs(1).Range
$startRange = $d.Range($firstPara.Start, $firstPara.End)
$d.Bookmarks.Add("_GoBack", $startRange)

# --- 2) Fix the "Heading 1" / "Titulo 1" style formatting --------------

$heading1 = $d.Styles("Ttulo1")
$heading1.ParagraphFormat.Alignment = 0   # wdAlignParagraphLeft -> <w:jc w:val="left"/>
